$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "2021"
$ws.Range("A46").ClearFormats()
$ws.Range("D46").Value = 52207.25
